$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: health
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 26
$ws.Range("D2").Value = 19.23076923076923

# Update row 3: police
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 66.66666666666666

# Insert a new row 4 for "fire" (this shifts the old row 4 "schools" down to row 5)
$ws.Rows.Item(4).Insert()

# New row 4: fire
$ws.Range("A4").Value = "fire"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 0

# Row 5: schools (shifted down), with updated values
$ws.Range("A5").Value = "schools"
$ws.Range("B5").Value = 23
$ws.Range("C5").Value = 138
$ws.Range("D5").Value = 16.66666666666666
